# Auto-generated Excel COM-interop script to apply market-data refresh
# Updates columns H, I, J (raw inputs) and recomputed K, L, M, N per row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2106.5881
$ws.Range("I62").Value = 1788.2858
$ws.Range("J62").Value = 3592
$ws.Range("K62").Value = 1788.2858
$ws.Range("L62").Value = 3592
$ws.Range("M62").Value = -1164.2858
$ws.Range("N62").Value = -4840
$ws.Range("H65").Value = 2106.5881
$ws.Range("I65").Value = 1788.2858
$ws.Range("J65").Value = 3592
$ws.Range("K65").Value = 8941.429
$ws.Range("L65").Value = 17960
$ws.Range("M65").Value = -5821.429
$ws.Range("N65").Value = -24200
$ws.Range("H103").Value = 1790.6
$ws.Range("J103").Value = 1851.5
$ws.Range("L103").Value = 5554.5
$ws.Range("N103").Value = -6726.5
$ws.Range("H112").Value = 1841.7646
$ws.Range("I112").Value = 2121.5
$ws.Range("K112").Value = 6364.5
$ws.Range("M112").Value = -5256.5
$ws.Range("H132").Value = 12480.4
$ws.Range("I132").Value = 13085.158
$ws.Range("K132").Value = 39255.474
$ws.Range("M132").Value = -36725.474
$ws.Range("H137").Value = 6255950.5
$ws.Range("I137").Value = 13335451
$ws.Range("J137").Value = 9332.294
$ws.Range("K137").Value = 40006353
$ws.Range("L137").Value = 27996.882
$ws.Range("M137").Value = -40003803
$ws.Range("N137").Value = -33096.882
$ws.Range("H138").Value = 6630.375
$ws.Range("J138").Value = 6339.39
$ws.Range("L138").Value = 19018.17
$ws.Range("N138").Value = -29298.17

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1114718.5
$ws.Range("I32").Value = 1191722.2
$ws.Range("J32").Value = 36666.668
$ws.Range("K32").Value = 1191722.2
$ws.Range("L32").Value = 36666.668
$ws.Range("M32").Value = -1191435.2
$ws.Range("N32").Value = -37240.668
$ws.Range("H45").Value = 2109.389
$ws.Range("I45").Value = 2147.375
$ws.Range("J45").Value = 1805.5
$ws.Range("K45").Value = 2147.375
$ws.Range("L45").Value = 1805.5
$ws.Range("M45").Value = -1770.375
$ws.Range("N45").Value = -2559.5
$ws.Range("H61").Value = 7956744.5
$ws.Range("I61").Value = 4331088.5
$ws.Range("K61").Value = 4331088.5
$ws.Range("M61").Value = -4330876.5
$ws.Range("H132").Value = 4017.1082
$ws.Range("I132").Value = 2490.48
$ws.Range("K132").Value = 7471.440000000001
$ws.Range("M132").Value = -4941.440000000001
$ws.Range("H136").Value = 7956744.5
$ws.Range("I136").Value = 4331088.5
$ws.Range("K136").Value = 12993265.5
$ws.Range("M136").Value = -12990715.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4087.875
$ws.Range("I105").Value = 2117.25
$ws.Range("K105").Value = 2117.25
$ws.Range("M105").Value = -370.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 408272.12
$ws.Range("I31").Value = 714608.2
$ws.Range("J31").Value = 3470.9285
$ws.Range("K31").Value = 714608.2
$ws.Range("L31").Value = 3470.9285
$ws.Range("M31").Value = -714313.2
$ws.Range("N31").Value = -4060.9285
$ws.Range("H34").Value = 408272.12
$ws.Range("I34").Value = 714608.2
$ws.Range("J34").Value = 3470.9285
$ws.Range("K34").Value = 714608.2
$ws.Range("L34").Value = 3470.9285
$ws.Range("M34").Value = -714406.2
$ws.Range("N34").Value = -3874.9285
$ws.Range("H58").Value = 6007790
$ws.Range("I58").Value = 18522580
$ws.Range("J58").Value = 1675747.8
$ws.Range("K58").Value = 18522580
$ws.Range("L58").Value = 1675747.8
$ws.Range("M58").Value = -18522377
$ws.Range("N58").Value = -1676153.8
$ws.Range("H122").Value = 3256978.2
$ws.Range("I122").Value = 3234.8696
$ws.Range("K122").Value = 9704.6088
$ws.Range("M122").Value = -7254.6088
$ws.Range("H136").Value = 6007790
$ws.Range("I136").Value = 18522580
$ws.Range("J136").Value = 1675747.8
$ws.Range("K136").Value = 55567740
$ws.Range("L136").Value = 5027243.4
$ws.Range("M136").Value = -55565190
$ws.Range("N136").Value = -5032343.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1849.5
$ws.Range("J52").Value = 1849.5
$ws.Range("L52").Value = 5548.5
$ws.Range("N52").Value = -6080.5
$ws.Range("H122").Value = 949942.9399999999
$ws.Range("I122").Value = 2016837.6
$ws.Range("J122").Value = 1592.1111
$ws.Range("K122").Value = 18151538.4
$ws.Range("L122").Value = 14328.9999
$ws.Range("M122").Value = -18149088.4
$ws.Range("N122").Value = -19228.9999
$ws.Range("H139").Value = 7266.5
$ws.Range("I139").Value = 3871.5386
$ws.Range("K139").Value = 11614.6158
$ws.Range("M139").Value = -6474.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3745.8333
$ws.Range("I102").Value = 2662.4285
$ws.Range("K102").Value = 2662.4285
$ws.Range("M102").Value = -1040.4285
$ws.Range("H132").Value = 9972.933999999999
$ws.Range("I132").Value = 7134.5483
$ws.Range("K132").Value = 21403.6449
$ws.Range("M132").Value = -18873.6449

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6099.0464
$ws.Range("I61").Value = 5552.7188
$ws.Range("K61").Value = 5552.7188
$ws.Range("M61").Value = -5350.7188
$ws.Range("H113").Value = 6099.0464
$ws.Range("I113").Value = 5552.7188
$ws.Range("K113").Value = 5552.7188
$ws.Range("M113").Value = -3382.7188
$ws.Range("H122").Value = 7551.0713
$ws.Range("I122").Value = 9906.143
$ws.Range("J122").Value = 5196
$ws.Range("K122").Value = 29718.429
$ws.Range("L122").Value = 15588
$ws.Range("M122").Value = -27268.429
$ws.Range("N122").Value = -20488
$ws.Range("H136").Value = 31977662
$ws.Range("I136").Value = 43627864
$ws.Range("J136").Value = 13337337
$ws.Range("K136").Value = 130883592
$ws.Range("L136").Value = 40012011
$ws.Range("M136").Value = -130881042
$ws.Range("N136").Value = -40017111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 27705
$ws.Range("J45").Value = 28787
$ws.Range("L45").Value = 28787
$ws.Range("N45").Value = -29769
$ws.Range("H54").Value = 516050
$ws.Range("J54").Value = 32100
$ws.Range("L54").Value = 32100
$ws.Range("N54").Value = -33140
$ws.Range("H136").Value = 3534640.8
$ws.Range("I136").Value = 3624880.2
$ws.Range("J136").Value = 3474481
$ws.Range("K136").Value = 10874640.6
$ws.Range("L136").Value = 10423443
$ws.Range("M136").Value = -10872090.6
$ws.Range("N136").Value = -10428543
